# This script applies a weekly data update to the "Apio" (celery) sheet.
# A new record is inserted at row 477 (pushing the existing rows 477-506
# down by one, to 478-507), and the new row 477 is populated with the
# latest week's price report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 506   # last row with data before the insert
$insertAt   = 477    # row at which the new record is inserted
$lastCol    = 18     # column R

# Shift rows [insertAt .. lastDataRow] down by one row, working from the
# bottom up so that we never overwrite a row before it has been read.
for ($r = $lastDataRow; $r -ge $insertAt; $r--) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($r + 1, $c)
        $v = $srcCell.Value()
        # Date cells come back as DateTime objects; convert back to the
        # underlying numeric serial and re-apply the source number format,
        # so the destination cell keeps the original date style instead of
        # Excel assigning (or omitting) one when the row is brand new.
        if ($v -is [DateTime]) {
            $dstCell.NumberFormat = $srcCell.NumberFormat
            $v = $v.ToOADate()
        }
        $dstCell.Value = $v
    }
}

# Populate the newly freed row with this week's record.
$ws.Cells.Item($insertAt, 1).Value  = 7
$ws.Cells.Item($insertAt, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($insertAt, 3).Value  = "Ñuble"
$ws.Cells.Item($insertAt, 4).Value  = 45267
$ws.Cells.Item($insertAt, 5).Value  = 16
$ws.Cells.Item($insertAt, 6).Value  = 100112017
$ws.Cells.Item($insertAt, 7).Value  = "Apio"
$ws.Cells.Item($insertAt, 8).Value  = "Americana (o)"
$ws.Cells.Item($insertAt, 9).Value  = "Segunda"
$ws.Cells.Item($insertAt, 10).Value = 200
$ws.Cells.Item($insertAt, 11).Value = 8000
$ws.Cells.Item($insertAt, 12).Value = 8000
$ws.Cells.Item($insertAt, 13).Value = 8000
$ws.Cells.Item($insertAt, 14).Value = "$/docena de matas"
$ws.Cells.Item($insertAt, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($insertAt, 16).Value = 1333
$ws.Cells.Item($insertAt, 17).Value = 6
$ws.Cells.Item($insertAt, 18).Value = "Hortaliza"

Write-Host "Inserted new row at $insertAt; rows $insertAt..$lastDataRow shifted to $($insertAt+1)..$($lastDataRow+1)"
